$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.592.49"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "1.752.88"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Formula = "'324.12"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Formula = "'1.001"
$ws.Range("D7").Formula = "'0.4633"
$ws.Range("E7").Value = "  +9.22%  "
$ws.Range("D8").Formula = "'0.3604"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Formula = "'0.07499"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Formula = "'42.16"
$ws.Range("E10").Value = "  -4.82%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Formula = "'1.001"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Formula = "'20.73"
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Formula = "'7.098"
$ws.Range("E15").Value = "  -3.41%  "
$ws.Range("D16").Value = "1.752.06"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").Formula = "'92.91"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Formula = "'16.77"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").Formula = "'5.826"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").Value = "27.632.36"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Formula = "'2.109"
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("D26").Formula = "'163.08"
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("D27").Formula = "'20.40"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").Value = "1.954.21"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").Formula = "'127.19"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").Formula = "'2.074"
$ws.Range("E30").Value = "  -2.90%  "
$ws.Range("D31").Formula = "'1.074"
$ws.Range("E31").Value = "  -8.48%  "
$ws.Range("D32").Formula = "'0.09247"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").Formula = "'3.671"
$ws.Range("E33").Value = "  +4.07%  "
$ws.Range("D34").Formula = "'5.523"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("E35").Value = "  -5.17%  "
$ws.Range("D36").Formula = "'0.02287"
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").Formula = "'0.2101"
$ws.Range("D38").Formula = "'0.06022"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").Formula = "'0.6350"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Formula = "'4.960"
$ws.Range("E40").Value = "  -1.58%  "
$ws.Range("D41").Formula = "'1.199"
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("D42").Formula = "'1.384"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").Formula = "'7.751"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").Formula = "'0.5883"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Formula = "'3.709"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Formula = "'122.49"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").Value = "  -3.45%  "
$ws.Range("D49").Formula = "'1.145"
$ws.Range("E49").Value = "  -4.23%  "
$ws.Range("D50").Formula = "'0.06854"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").Formula = "'72.09"
$ws.Range("E51").Value = "  -2.79%  "
